$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.7314747453968576

$ws.Range("C3").Value = "cornstover"
$ws.Range("D3").Value = 0.7275220018020688
$ws.Range("E3").Value = 0.5805274874239775
$ws.Range("F3").Value = 0.7067329609574093
$ws.Range("G3").Value = 0.8432608888730624

$ws.Range("C4").Value = "sugarcane"
$ws.Range("D4").Value = 0.729124654736764
$ws.Range("E4").Value = 0.5950268706740716
$ws.Range("F4").Value = 0.718866082424461
$ws.Range("G4").Value = 0.8570247350279394

$ws.Range("D5").Value = 0.7306849292360015

$ws.Range("D6").Value = 0.7199420883218625

$ws.Range("C7").Value = "cornstover"
$ws.Range("D7").Value = 0.7160516452034039
$ws.Range("E7").Value = 0.5551076458346899
$ws.Range("F7").Value = 0.6949337710892807
$ws.Range("G7").Value = 0.8406469372339114

$ws.Range("C8").Value = "sugarcane"
$ws.Range("D8").Value = 0.7170303666175638
$ws.Range("E8").Value = 0.5390702079497377
$ws.Range("F8").Value = 0.6767208248177626
$ws.Range("G8").Value = 0.8279481166749891

$ws.Range("D9").Value = 0.7191651251354755

$ws.Range("D10").Value = 0.6031912016243293

$ws.Range("C11").Value = "cornstover"
$ws.Range("D11").Value = 0.6028488348905511
$ws.Range("E11").Value = 0.5878796681521615
$ws.Range("F11").Value = 0.5885187862166137
$ws.Range("G11").Value = 0.5891170179878673

$ws.Range("C12").Value = "sugarcane"
$ws.Range("D12").Value = 0.6027367034234429
$ws.Range("E12").Value = 0.5938542188089596
$ws.Range("F12").Value = 0.5939109400288758
$ws.Range("G12").Value = 0.5939638334198035

$ws.Range("D13").Value = 0.6032564221695031

$ws.Range("D14").Value = 0.626369466760216

$ws.Range("C15").Value = "cornstover"
$ws.Range("D15").Value = 0.6258257611068714
$ws.Range("E15").Value = 0.6114213630217711
$ws.Range("F15").Value = 0.612007871668953
$ws.Range("G15").Value = 0.6126434268673036

$ws.Range("C16").Value = "sugarcane"
$ws.Range("D16").Value = 0.6258731801449157
$ws.Range("E16").Value = 0.6164189239676208
$ws.Range("F16").Value = 0.6165417412970111
$ws.Range("G16").Value = 0.6170759411149578

$ws.Range("D17").Value = 0.6263624387182642

$ws.Range("D18").Value = 0.6810615160036566
$ws.Range("E18").Value = 0.6737624744204248
$ws.Range("F18").Value = 0.6743112509176761
$ws.Range("G18").Value = 0.6745486377266964

$ws.Range("C19").Value = "cornstover"
$ws.Range("D19").Value = 0.6937248387038393
$ws.Range("E19").Value = 0.6719516411699592
$ws.Range("F19").Value = 0.6721268095177348
$ws.Range("G19").Value = 0.6737795309311089

$ws.Range("C20").Value = "sugarcane"
$ws.Range("D20").Value = 0.6796777814820135
$ws.Range("E20").Value = 0.6533136550058209
$ws.Range("F20").Value = 0.6535833496110711
$ws.Range("G20").Value = 0.6538453745643791

$ws.Range("D21").Value = 0.6841117478261509
$ws.Range("E21").Value = 0.6714429719355082
$ws.Range("F21").Value = 0.6716585489842333
$ws.Range("G21").Value = 0.6718791796141838

$ws.Range("D22").Value = 0.7070817008827719
$ws.Range("E22").Value = 0.6979255606208118
$ws.Range("F22").Value = 0.6981345699727914
$ws.Range("G22").Value = 0.6991365802059701

$ws.Range("C23").Value = "cornstover"
$ws.Range("D23").Value = 0.7150761631438003
$ws.Range("E23").Value = 0.6907572379493958
$ws.Range("F23").Value = 0.6920983812852993
$ws.Range("G23").Value = 0.6934675501331323

$ws.Range("C24").Value = "sugarcane"
$ws.Range("D24").Value = 0.7056682221142379
$ws.Range("E24").Value = 0.6848391174624211
$ws.Range("F24").Value = 0.685441725242266
$ws.Range("G24").Value = 0.6868155028790894

$ws.Range("D25").Value = 0.7086923329599475
$ws.Range("E25").Value = 0.6882993284376934
$ws.Range("F25").Value = 0.688451679968975
$ws.Range("G25").Value = 0.6886189202870117
